$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 45311.3
$ws.Range("I6").Value = 42.5
$ws.Range("J6").Value = 75490.5
$ws.Range("K6").Value = 127.5
$ws.Range("L6").Value = 226471.5
$ws.Range("M6").Value = -15.5
$ws.Range("N6").Value = -226695.5

$ws.Range("H7").Value = 9250
$ws.Range("J7").Value = 9250
$ws.Range("L7").Value = 9250
$ws.Range("N7").Value = -9474

$ws.Range("H9").Value = 200.2
$ws.Range("I9").Value = 175.25
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 175.25
$ws.Range("L9").Value = 300
$ws.Range("M9").Value = -6.25
$ws.Range("N9").Value = -638

$ws.Range("H12").Value = 593.3333
$ws.Range("I12").Value = 540
$ws.Range("J12").Value = 700
$ws.Range("K12").Value = 540
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = -370
$ws.Range("N12").Value = -1040

$ws.Range("H14").Value = 9250
$ws.Range("J14").Value = 9250
$ws.Range("L14").Value = 9250
$ws.Range("N14").Value = -9632

$ws.Range("H21").Value = 80019
$ws.Range("I21").Value = 80019
$ws.Range("K21").Value = 80019
$ws.Range("M21").Value = -79551

$ws.Range("H23").Value = 80019
$ws.Range("I23").Value = 80019
$ws.Range("K23").Value = 80019
$ws.Range("M23").Value = -79785

$ws.Range("H29").Value = 1680
$ws.Range("I29").Value = 133.33333
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 399.99999
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = -118.99999
$ws.Range("N29").Value = -12562

$ws.Range("H38").Value = 871.95
$ws.Range("I38").Value = 110
$ws.Range("J38").Value = 2287
$ws.Range("K38").Value = 330
$ws.Range("L38").Value = 6861
$ws.Range("M38").Value = 42
$ws.Range("N38").Value = -7605

$ws.Range("H58").Value = 274
$ws.Range("I58").Value = 274
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 822
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -672

$ws.Range("H86").Value = 2764.75
$ws.Range("I86").Value = 2254.1667
$ws.Range("J86").Value = 3071.1
$ws.Range("K86").Value = 2254.1667
$ws.Range("L86").Value = 3071.1
$ws.Range("M86").Value = -1131.1667
$ws.Range("N86").Value = -5317.1

$ws.Range("H87").Value = 39630
$ws.Range("J87").Value = 39630
$ws.Range("L87").Value = 39630
$ws.Range("N87").Value = -42126

$ws.Range("H89").Value = 2764.75
$ws.Range("I89").Value = 2254.1667
$ws.Range("J89").Value = 3071.1
$ws.Range("K89").Value = 11270.8335
$ws.Range("L89").Value = 15355.5
$ws.Range("M89").Value = -5654.833500000001
$ws.Range("N89").Value = -26587.5

$ws.Range("H90").Value = 39630
$ws.Range("J90").Value = 39630
$ws.Range("L90").Value = 118890
$ws.Range("N90").Value = -131370

$ws.Range("H124").Value = 34705.883
$ws.Range("J124").Value = 34705.883
$ws.Range("L124").Value = 34705.883
$ws.Range("N124").Value = -44525.883

$ws.Range("H131").Value = 1257.9736
$ws.Range("I131").Value = 1023.7143
$ws.Range("J131").Value = 1394.625
$ws.Range("K131").Value = 3071.1429
$ws.Range("L131").Value = 4183.875
$ws.Range("M131").Value = 1968.8571
$ws.Range("N131").Value = -14263.875

$ws.Range("H138").Value = 3512271
$ws.Range("I138").Value = 2870.5
$ws.Range("J138").Value = 4881793.5
$ws.Range("K138").Value = 8611.5
$ws.Range("L138").Value = 14645380.5
$ws.Range("M138").Value = -3471.5
$ws.Range("N138").Value = -14655660.5

$ws.Range("H141").Value = 2715.9033
$ws.Range("I141").Value = 2578.3215
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 7734.9645
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -2554.9645
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13483.31
$ws.Range("I32").Value = 11273.823
$ws.Range("K32").Value = 11273.823
$ws.Range("M32").Value = -10986.823

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 48880.5
$ws.Range("J88").Value = 48880.5
$ws.Range("L88").Value = 48880.5
$ws.Range("N88").Value = -49692.5

$ws.Range("H91").Value = 48880.5
$ws.Range("J91").Value = 48880.5
$ws.Range("L91").Value = 48880.5
$ws.Range("N91").Value = -51688.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2975.375
$ws.Range("I86").Value = 1998.75
$ws.Range("J86").Value = 3952
$ws.Range("K86").Value = 1998.75
$ws.Range("L86").Value = 3952
$ws.Range("M86").Value = -875.75
$ws.Range("N86").Value = -6198

$ws.Range("H89").Value = 2975.375
$ws.Range("I89").Value = 1998.75
$ws.Range("J89").Value = 3952
$ws.Range("K89").Value = 9993.75
$ws.Range("L89").Value = 19760
$ws.Range("M89").Value = -4377.75
$ws.Range("N89").Value = -30992

$ws.Range("H129").Value = 30444.875
$ws.Range("J129").Value = 30444.875
$ws.Range("L129").Value = 30444.875
$ws.Range("N129").Value = -40444.875

$ws.Range("H134").Value = 72742.47
$ws.Range("I134").Value = 562.4286
$ws.Range("J134").Value = 135900
$ws.Range("K134").Value = 1687.2858
$ws.Range("L134").Value = 407700
$ws.Range("M134").Value = 847.7142000000001
$ws.Range("N134").Value = -412770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1013.6923
$ws.Range("I9").Value = 563
$ws.Range("J9").Value = 1295.375
$ws.Range("K9").Value = 563
$ws.Range("L9").Value = 1295.375
$ws.Range("M9").Value = -339
$ws.Range("N9").Value = -1743.375

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H136").Value = 38120.414
$ws.Range("I136").Value = 22865.076
$ws.Range("K136").Value = 68595.228
$ws.Range("M136").Value = -66045.228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16040

$ws.Range("H81").Value = 2164.2856
$ws.Range("I81").Value = 700
$ws.Range("J81").Value = 2750
$ws.Range("K81").Value = 1400
$ws.Range("L81").Value = 5500
$ws.Range("M81").Value = -339
$ws.Range("N81").Value = -7622

$ws.Range("H84").Value = 2164.2856
$ws.Range("I84").Value = 700
$ws.Range("J84").Value = 2750
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 27500
$ws.Range("M84").Value = -1696
$ws.Range("N84").Value = -38108

$ws.Range("H132").Value = 43065.832
$ws.Range("I132").Value = 31292.758
$ws.Range("J132").Value = 68966.60000000001
$ws.Range("K132").Value = 93878.274
$ws.Range("L132").Value = 206899.8
$ws.Range("M132").Value = -91348.274
$ws.Range("N132").Value = -211959.8
